$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect the sheet first so locked cells (footer text, Weight/% Change columns) can be edited
$ws.Unprotect("D382")

# Update the confidentiality footer date
$footerCell = $ws.Range("A44")
$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."
if ($footerCell.Text -eq $oldText) {
    $footerCell.Value = $newText
}

# Update Weight (D) and Percent Change (E) figures for each holding row
$ws.Range("D2").Value = 0.07774594805516165
$ws.Range("E2").Value = 0.005141388174807249
$ws.Range("D3").Value = 0.06776333381737981
$ws.Range("E3").Value = 0.0178187139536885
$ws.Range("D4").Value = 0.05388726670213612
$ws.Range("E4").Value = 0.02400797607178462
$ws.Range("D5").Value = 0.04676103833187557
$ws.Range("E5").Value = 0.001894125911773736
$ws.Range("D6").Value = 0.04449210886563351
$ws.Range("E6").Value = 0.01665027859718116
$ws.Range("D7").Value = 0.03830885009899415
$ws.Range("E7").Value = 0.02673462518632497
$ws.Range("D8").Value = 0.03978063791594354
$ws.Range("E8").Value = 0.01827498919553006
$ws.Range("D9").Value = 0.03497061796974945
$ws.Range("E9").Value = 0.008357585254831656
$ws.Range("D10").Value = 0.03119889305837447
$ws.Range("E10").Value = -0.00380327288113913
$ws.Range("D11").Value = 0.02869126625636309
$ws.Range("E11").Value = -0.005296666845005626
$ws.Range("D12").Value = 0.03121700642533317
$ws.Range("E12").Value = 0.027084439723845
$ws.Range("D13").Value = 0.02995121978176156
$ws.Range("E13").Value = 0.0153548109349213
$ws.Range("D14").Value = 0.02750468721745834
$ws.Range("E14").Value = 0.02171001227815594
$ws.Range("D15").Value = 0.03073439264806065
$ws.Range("E15").Value = 0.01475876535810605
$ws.Range("D16").Value = 0.02848173451145946
$ws.Range("E16").Value = 0.01452477835565502
$ws.Range("D17").Value = 0.02737083250569573
$ws.Range("E17").Value = 0.02230971128608927
$ws.Range("D18").Value = 0.02284863089041321
$ws.Range("E18").Value = 0.04595291841341509
$ws.Range("D19").Value = 0.01966037129946216
$ws.Range("E19").Value = -0.03385436999328528
$ws.Range("D20").Value = 0.02333400772366247
$ws.Range("E20").Value = -0.01989342806394312
$ws.Range("D21").Value = 0.02155951177414903
$ws.Range("E21").Value = 0.02723389106443563
$ws.Range("D22").Value = 0.02202002110360754
$ws.Range("E22").Value = 0.01394214011850825
$ws.Range("D23").Value = 0.02076160260930726
$ws.Range("E23").Value = 0.01960784313725483
$ws.Range("D24").Value = 0.02031644359083071
$ws.Range("E24").Value = 0.007646276595744572
$ws.Range("D25").Value = 0.01827454522401169
$ws.Range("E25").Value = 0.02536749265014682
$ws.Range("D26").Value = 0.01777197604246264
$ws.Range("E26").Value = 0.0338152985074629
$ws.Range("D27").Value = 0.01938360519245759
$ws.Range("E27").Value = 0.0407840031676896
$ws.Range("D28").Value = 0.01774142892360856
$ws.Range("E28").Value = 0.02293708955934148
$ws.Range("D29").Value = 0.0186831705023512
$ws.Range("E29").Value = 0.01735243854344692
$ws.Range("D30").Value = 0.01806900455996338
$ws.Range("E30").Value = 0.01626016260162588
$ws.Range("D31").Value = 0.0179093613257511
$ws.Range("E31").Value = 0.02360483753460585
$ws.Range("D32").Value = 0.01594590304804985
$ws.Range("E32").Value = 0.02021563342318067
$ws.Range("D33").Value = 0.01714691137927763
$ws.Range("E33").Value = 0.02198667908042684
$ws.Range("D34").Value = 0.007696799429459641
$ws.Range("E34").Value = 0.0242516104585071
$ws.Range("D35").Value = 0.007719057380383469
$ws.Range("E35").Value = 0.01032096408543137
$ws.Range("D36").Value = 0.007191313688824022
$ws.Range("E36").Value = 0.02898736338797803
$ws.Range("D37").Value = 0.006307135776263691
$ws.Range("E37").Value = 0.01776674454828653
$ws.Range("D38").Value = 0.007087238580366399
$ws.Range("E38").Value = 0.04405458089668612
$ws.Range("D39").Value = 0.006922836749749713
$ws.Range("E39").Value = 0.04013392758154266
$ws.Range("D40").Value = 0.006789289044206746
$ws.Range("E40").Value = 0.0401772592642835
$ws.Range("D41").Value = 0.9999999999999999
$ws.Range("E41").Value = 0.01571211781179471

# Restore sheet protection with the original password
$ws.Protect("D382")
